# Fruta / hortaliza, semanal
#
# A new week's worth of price observations (date 44694) was inserted for
# "Especial" and "Primera" quality grades, ahead of the existing data block
# that starts at row 610. This pushes the existing rows (610 downward)
# down by two rows (new dimension A1:T728, was A1:T726).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 610:611, shifting rows 610..726 down to 612..728
$ws.Rows("610:611").Insert()

# --- New row 610: Especial, fecha 44694 ---
$ws.Range("A610").Value = 8
$ws.Range("B610").Value = "Terminal La Palmera de La Serena"
$ws.Range("C610").Value = "Coquimbo"
$ws.Range("D610").Value = 44694
$ws.Range("D610").NumberFormat = $ws.Range("D612").NumberFormat
$ws.Range("E610").Value = 4
$ws.Range("F610").Value = "Fruta"
$ws.Range("G610").Value = 100101
$ws.Range("H610").Value = "Berries"
$ws.Range("I610").Value = 100112025
$ws.Range("J610").Value = "Frutilla"
$ws.Range("K610").Value = "Sin especificar"
$ws.Range("L610").Value = "Especial"
$ws.Range("M610").Value = 400
$ws.Range("N610").Value = 18000
$ws.Range("O610").Value = 19000
$ws.Range("P610").Value = 18500
$ws.Range("Q610").Value = "$/bandeja 7 kilos"
$ws.Range("R610").Value = "Provincia de Melipilla"
$ws.Range("S610").Value = 2643
$ws.Range("T610").Value = 7

# --- New row 611: Primera, fecha 44694 ---
$ws.Range("A611").Value = 8
$ws.Range("B611").Value = "Terminal La Palmera de La Serena"
$ws.Range("C611").Value = "Coquimbo"
$ws.Range("D611").Value = 44694
$ws.Range("D611").NumberFormat = $ws.Range("D612").NumberFormat
$ws.Range("E611").Value = 4
$ws.Range("F611").Value = "Fruta"
$ws.Range("G611").Value = 100101
$ws.Range("H611").Value = "Berries"
$ws.Range("I611").Value = 100112025
$ws.Range("J611").Value = "Frutilla"
$ws.Range("K611").Value = "Sin especificar"
$ws.Range("L611").Value = "Primera"
$ws.Range("M611").Value = 500
$ws.Range("N611").Value = 15000
$ws.Range("O611").Value = 16000
$ws.Range("P611").Value = 15500
$ws.Range("Q611").Value = "$/bandeja 7 kilos"
$ws.Range("R611").Value = "Provincia de Melipilla"
$ws.Range("S611").Value = 2214
$ws.Range("T611").Value = 7
